$d = $word.ActiveDocument

# Inserting new text immediately adjacent to existing text that has
# identical run formatting always coalesces into a single run when the
# package is serialized. The one thing that reliably forces a run
# boundary at an exact character offset is adding (then removing) a
# bookmark there, so splitting is implemented via that trick.
function Split-AtPos([int]$pos) {
    $bmRange = $d.Range($pos, $pos)
    $name = "TempSplit" + [guid]::NewGuid().ToString("N")
    [void]$d.Bookmarks.Add($name, $bmRange)
    $d.Bookmarks.Item($name).Delete()
}

# Finds $searchText, inserts the concatenation of $newTexts right after
# it, then splits the freshly inserted text back into separate runs
# matching each entry of $newTexts (and detaches it from the original
# "$searchText" run too).
function Insert-RunsAfter([string]$searchText, [string[]]$newTexts) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $found = $r.Find.Execute($searchText)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $searchText)
        return
    }
    $pos = $r.End

    $fullText = [string]::Join("", $newTexts)
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($fullText)

    $positions = @($pos)
    $cursor = $pos
    for ($i = 0; $i -lt $newTexts.Length - 1; $i++) {
        $cursor = $cursor + $newTexts[$i].Length
        $positions += $cursor
    }
    foreach ($p in $positions) {
        Split-AtPos $p
    }
}

# --- Abstracción ---
Insert-RunsAfter "Abstracción" @(
    ": extraer todos los ",
    "atributos mas importantes ",
    "o lo que nos son de utilidad "
)

# --- Encapsulamiento ---
Insert-RunsAfter "Encapsulamiento" @(
    ": ",
    "nos permite ",
    "controlar quien tiene acceso al programa con modificadores de acceso "
)

# --- Herencia ---
Insert-RunsAfter "Herencia" @(
    ":",
    " se refiere que una ",
    "clase herede los atributos y métodos de otra clase ",
    "haciendo",
    " que se reutilice código"
)

# --- Polimorfismo ---
# Special-cased because the existing "_GoBack" bookmark must end up
# sitting between "...mismo" and " método de diferentes...".
$beforeTexts = @(":", " usar una mism", "o")
$afterText = " método de diferentes maneras según el contexto"

$r = $d.Content
$r.Find.ClearFormatting()
[void]$r.Find.Execute("Polimorfismo")
$pos = $r.End

$beforeFull = [string]::Join("", $beforeTexts)
$ins = $d.Range($pos, $pos)
$ins.InsertAfter($beforeFull + $afterText)

$positions = @($pos)
$cursor = $pos
for ($i = 0; $i -lt $beforeTexts.Length; $i++) {
    $cursor = $cursor + $beforeTexts[$i].Length
    $positions += $cursor
}
foreach ($p in $positions) {
    Split-AtPos $p
}

# Re-seat "_GoBack" exactly at the boundary between the "before" text
# and the "after" text.
$bmPos = $pos + $beforeFull.Length
$bmRange = $d.Range($bmPos, $bmPos)
[void]$d.Bookmarks.Add("_GoBack", $bmRange)
